# Update "report_co_so" workbook:
#   Sheet 1 (DOANH SỐ CÁ NHÂN): insert two new columns ("Công phụ phẫu 1" at G,
#     "Công phụ phẫu 2" at I - pushing the old "Số lần phụ phẫu 2" to H and
#     "Doanh số thu nợ" to J), populate the new columns, tweak a couple of
#     existing values, and append a "Tổng" (Total) row.
#   Sheet 2 (CHI TIÊU): update a few expense figures.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: insert the two new columns -----------------------------------
# Inserting at column 7 (G) pushes old G -> H and old H -> I.
$ws1.Columns.Item(7).Insert()
# Inserting at column 9 (I) pushes (the just-shifted) old H -> J, leaving a
# fresh, empty column I for "Công phụ phẫu 2".
$ws1.Columns.Item(9).Insert()

# --- Sheet 1: headers --------------------------------------------------------
$ws1.Range("G1").Value = "Công phụ phẫu 1"
$ws1.Range("I1").Value = "Công phụ phẫu 2"

# --- Sheet 1: new column G ("Công phụ phẫu 1") values per row --------------
$ws1.Range("G2").Value = 0
$ws1.Range("G3").Value = 0
$ws1.Range("G4").Value = 900000
$ws1.Range("G5").Value = 1900000
$ws1.Range("G6").Value = 0
$ws1.Range("G7").Value = 0
$ws1.Range("G8").Value = 0
$ws1.Range("G9").Value = 0
$ws1.Range("G10").Value = 1350000
$ws1.Range("G11").Value = 0
$ws1.Range("G12").Value = 0
$ws1.Range("G13").Value = 0
$ws1.Range("G14").Value = 50000
$ws1.Range("G15").Value = 0
$ws1.Range("G16").Value = 550000
$ws1.Range("G17").Value = 0
$ws1.Range("G18").Value = 0

# --- Sheet 1: new column I ("Công phụ phẫu 2") values per row --------------
$ws1.Range("I2").Value = 0
$ws1.Range("I3").Value = 0
$ws1.Range("I4").Value = 0
$ws1.Range("I5").Value = 100000
$ws1.Range("I6").Value = 0
$ws1.Range("I7").Value = 0
$ws1.Range("I8").Value = 0
$ws1.Range("I9").Value = 0
$ws1.Range("I10").Value = 0
$ws1.Range("I11").Value = 50000
$ws1.Range("I12").Value = 0
$ws1.Range("I13").Value = 0
$ws1.Range("I14").Value = 0
$ws1.Range("I15").Value = 100000
$ws1.Range("I16").Value = 0
$ws1.Range("I17").Value = 0
$ws1.Range("I18").Value = 0

# --- Sheet 1: other value corrections ---------------------------------------
$ws1.Range("D3").Value = 74300000
$ws1.Range("D17").Value = 36100000

# --- Sheet 1: append the new "Tổng" (Total) row 19 --------------------------
$ws1.Range("A19").Value = "Tổng"
$ws1.Range("B19").Value = 1043200000
$ws1.Range("C19").Value = 0
$ws1.Range("D19").Value = 868200000
$ws1.Range("E19").Value = 324000000
$ws1.Range("F19").Value = 61
$ws1.Range("G19").Value = 4750000
$ws1.Range("H19").Value = 8
$ws1.Range("I19").Value = 250000
$ws1.Range("J19").Value = 143900000

# --- Sheet 2: expense figure corrections ------------------------------------
$ws2.Range("B3").Value = 29432000
$ws2.Range("B10").Value = 54242000
$ws2.Range("B12").Value = 331227000
